# "Re-did some comet times for consistency"
# Updates a handful of raw timing measurements on the "Comet" sheet,
# lets the dependent Speedup/Efficiency formulas recalc on their own,
# normalizes the number format on the touched columns (to the
# "0.0000" / "0.000000" formats already used elsewhere in the
# workbook), and moves the active-tab/selection from "Linux" back to
# "Comet".

$wb = $excel.ActiveWorkbook
$comet = $wb.Worksheets.Item("Comet")
$linux = $wb.Worksheets.Item("Linux")

# ---------------------------------------------------------------
# 1. Re-measured raw times (columns are matrix sizes, rows are
#    thread counts 1/2/4/8/12/16 -> sheet rows 4..9).
# ---------------------------------------------------------------
$comet.Range("I5").Value = 0.0726
$comet.Range("G6").Value = 0.011
$comet.Range("J7").Value = 0.1152
$comet.Range("I8").Value = 0.0628
$comet.Range("J8").Value = 0.1105
$comet.Range("G9").Value = 0.0134
$comet.Range("I9").Value = 0.0597
$comet.Range("J9").Value = 0.1086

# ---------------------------------------------------------------
# 2. Number-format touch-ups on column C (and a few other cells
#    incidentally re-touched alongside their value) in the raw
#    "Times" block - matches the "0.0000" format already used on
#    other columns in this workbook.
# ---------------------------------------------------------------
$comet.Range("C4").NumberFormat = "0.0000"
$comet.Range("J4").NumberFormat = "0.0000"
$comet.Range("C5").NumberFormat = "0.0000"
$comet.Range("F5").NumberFormat = "0.0000"
$comet.Range("G6").NumberFormat = "0.0000"
$comet.Range("J7").NumberFormat = "0.0000"
$comet.Range("C8").NumberFormat = "0.0000"
$comet.Range("G8").NumberFormat = "0.0000"
$comet.Range("C9").NumberFormat = "0.0000"
$comet.Range("J9").NumberFormat = "0.0000"

# Same touch-up (using the finer "0.000000" format) on the computed
# Speedup / Efficiency blocks that derive from column C (and the
# couple of other cells whose format drifted alongside theirs).
$comet.Range("C15").NumberFormat = "0.000000"
$comet.Range("C16").NumberFormat = "0.000000"
$comet.Range("C18").NumberFormat = "0.000000"
$comet.Range("E18").NumberFormat = "0.000000"
$comet.Range("C19").NumberFormat = "0.000000"
$comet.Range("G19").NumberFormat = "0.000000"
$comet.Range("I19").NumberFormat = "0.000000"
$comet.Range("C25").NumberFormat = "0.000000"
$comet.Range("C26").NumberFormat = "0.000000"
$comet.Range("C28").NumberFormat = "0.000000"
$comet.Range("C29").NumberFormat = "0.000000"

# ---------------------------------------------------------------
# 3. Move the active tab / selection back to "Comet" (it had been
#    left on "Linux").
# ---------------------------------------------------------------
$comet.Select()
$comet.Range("U8").Select()
